$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$para1 = $tr.Paragraphs(1,1)
$para1.Text = "Profesorju Gregorju Medetu"
